$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("H1").Value = "Avg_Experiment_Time"
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"
$ws.Range("O1").Value = "Obs_Prob"

# Apply header style (bold, centered, bordered) to the newly added header cells I1:O1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows 2-13 ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 55.298
$ws.Range("D2").Value = 55.298
$ws.Range("E2").Value = 3.0855516
$ws.Range("F2").Value = 0.15053798
$ws.Range("G2").Value = 0.15053798
$ws.Range("H2").Value = 8.307759420000002
$ws.Range("I2").Value = 5.236031708852824
$ws.Range("J2").Value = 5.236031708852824
$ws.Range("K2").Value = 0.3156637658242528
$ws.Range("L2").Value = 0.01905519671804852
$ws.Range("M2").Value = 0.01905519671804852
$ws.Range("N2").Value = 1.185786966129223
$ws.Range("O2").Value = 0.15

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 86.22
$ws.Range("D3").Value = 86.22
$ws.Range("E3").Value = 1.99355016
$ws.Range("F3").Value = 0.09323035999999998
$ws.Range("G3").Value = 0.09323035999999998
$ws.Range("H3").Value = 7.91036976
$ws.Range("I3").Value = 11.36212108680339
$ws.Range("J3").Value = 11.36212108680339
$ws.Range("K3").Value = 0.2577630931431959
$ws.Range("L3").Value = 0.01437619874832365
$ws.Range("M3").Value = 0.01437619874832365
$ws.Range("N3").Value = 0.7258704967389256
$ws.Range("O3").Value = 0.85

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 30.84
$ws.Range("D4").Value = 61.654
$ws.Range("E4").Value = 2.798544619999999
$ws.Range("F4").Value = 0.2196498
$ws.Range("G4").Value = 0.10982498
$ws.Range("H4").Value = 3.3473545
$ws.Range("I4").Value = 4.786671541781177
$ws.Range("J4").Value = 9.570028992109672
$ws.Range("K4").Value = 0.3789214474499044
$ws.Range("L4").Value = 0.03084583071510654
$ws.Range("M4").Value = 0.01542299750059653
$ws.Range("N4").Value = 0.4575747084627816
$ws.Range("O4").Value = 0.15

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 46.318
$ws.Range("D5").Value = 90.794
$ws.Range("E5").Value = 1.9068
$ws.Range("F5").Value = 0.14300214
$ws.Range("G5").Value = 0.0715008
$ws.Range("H5").Value = 3.23749456
$ws.Range("I5").Value = 7.87533604588854
$ws.Range("J5").Value = 14.27189328280419
$ws.Range("K5").Value = 0.2958295109652109
$ws.Range("L5").Value = 0.02432820726971711
$ws.Range("M5").Value = 0.01216412596453037
$ws.Range("N5").Value = 0.3465463592563541
$ws.Range("O5").Value = 0.85

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 15.236
$ws.Range("D6").Value = 60.862
$ws.Range("E6").Value = 2.87115312
$ws.Range("F6").Value = 0.29655032
$ws.Range("G6").Value = 0.07413755999999999
$ws.Range("H6").Value = 1.11953054
$ws.Range("I6").Value = 2.97422191370852
$ws.Range("J6").Value = 11.86898063874087
$ws.Range("K6").Value = 0.4997529037485183
$ws.Range("L6").Value = 0.05562049348656463
$ws.Range("M6").Value = 0.01390513878447933
$ws.Range("N6").Value = 0.2632440558359758
$ws.Range("O6").Value = 0.15

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 23.976
$ws.Range("D7").Value = 89.018
$ws.Range("E7").Value = 1.95328032
$ws.Range("F7").Value = 0.20806374
$ws.Range("G7").Value = 0.0520158
$ws.Range("H7").Value = 1.20999708
$ws.Range("I7").Value = 5.13577417267401
$ws.Range("J7").Value = 15.28649829294205
$ws.Range("K7").Value = 0.3292600327434565
$ws.Range("L7").Value = 0.04182342545567128
$ws.Range("M7").Value = 0.01045617936599539
$ws.Range("N7").Value = 0.1906400623799785
$ws.Range("O7").Value = 0.85

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 10.04
$ws.Range("D8").Value = 60.088
$ws.Range("E8").Value = 2.89505116
$ws.Range("F8").Value = 0.33340368
$ws.Range("G8").Value = 0.05556754
$ws.Range("H8").Value = 0.5627143999999999
$ws.Range("I8").Value = 1.756294149132358
$ws.Range("J8").Value = 10.46629555694711
$ws.Range("K8").Value = 0.4951676843845648
$ws.Range("L8").Value = 0.07580552601761441
$ws.Range("M8").Value = 0.01263413292545064
$ws.Range("N8").Value = 0.1835178211662595
$ws.Range("O8").Value = 0.15

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 16.786
$ws.Range("D9").Value = 87.16
$ws.Range("E9").Value = 1.99036614
$ws.Range("F9").Value = 0.22363896
$ws.Range("G9").Value = 0.03727318
$ws.Range("H9").Value = 0.61069242
$ws.Range("I9").Value = 4.041394529055404
$ws.Range("J9").Value = 14.34985144428143
$ws.Range("K9").Value = 0.3209049513642951
$ws.Range("L9").Value = 0.0454976311933579
$ws.Range("M9").Value = 0.007582819924250292
$ws.Range("N9").Value = 0.134390210815631
$ws.Range("O9").Value = 0.85

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 7.462
$ws.Range("D10").Value = 59.526
$ws.Range("E10").Value = 2.94189828
$ws.Range("F10").Value = 0.30860202
$ws.Range("G10").Value = 0.0385753
$ws.Range("H10").Value = 0.2948432
$ws.Range("I10").Value = 1.427105202329889
$ws.Range("J10").Value = 11.35990399773147
$ws.Range("K10").Value = 0.559584902010909
$ws.Range("L10").Value = 0.0917375751385669
$ws.Range("M10").Value = 0.01146720800495298
$ws.Range("N10").Value = 0.1314510021860291
$ws.Range("O10").Value = 0.15

# Row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 12.944
$ws.Range("D11").Value = 81.996
$ws.Range("E11").Value = 2.11643684
$ws.Range("F11").Value = 0.21395906
$ws.Range("G11").Value = 0.02674488
$ws.Range("H11").Value = 0.33834146
$ws.Range("I11").Value = 3.339336211384731
$ws.Range("J11").Value = 13.48330757515873
$ws.Range("K11").Value = 0.3456763062377906
$ws.Range("L11").Value = 0.05146242290895095
$ws.Range("M11").Value = 0.006432552374111803
$ws.Range("N11").Value = 0.09306303523484846
$ws.Range("O11").Value = 0.85

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 6.166
$ws.Range("D12").Value = 61.474
$ws.Range("E12").Value = 2.86728752
$ws.Range("F12").Value = 0.2839774999999999
$ws.Range("G12").Value = 0.02839776
$ws.Range("H12").Value = 0.18065502
$ws.Range("I12").Value = 1.279734128199196
$ws.Range("J12").Value = 12.77025774860993
$ws.Range("K12").Value = 0.5937619162538437
$ws.Range("L12").Value = 0.09080514557267638
$ws.Range("M12").Value = 0.009080204425408848
$ws.Range("N12").Value = 0.08488609587172853
$ws.Range("O12").Value = 0.15

# Row 13
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 10.846
$ws.Range("D13").Value = 78.60599999999999
$ws.Range("E13").Value = 2.21957778
$ws.Range("F13").Value = 0.18920812
$ws.Range("G13").Value = 0.01892072
$ws.Range("H13").Value = 0.20157828
$ws.Range("I13").Value = 2.914665426613867
$ws.Range("J13").Value = 14.14407557062302
$ws.Range("K13").Value = 0.4002347793078543
$ws.Range("L13").Value = 0.05029513277870987
$ws.Range("M13").Value = 0.005029555122579145
$ws.Range("N13").Value = 0.06669907373506646
$ws.Range("O13").Value = 0.85

